$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.348.99"
$ws.Range("E2").Value = "  -2.31%  "
$ws.Range("D3").Value = "2.437.59"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").Value = "580.76"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").Value = "143.00"
$ws.Range("E6").Value = "  -3.88%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -2.90%  "
$ws.Range("D9").Value = "2.434.92"
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("E10").Value = "  -5.15%  "
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").Value = "5.19"
$ws.Range("E12").Value = "  -1.74%  "
$ws.Range("D13").Value = "0.345"
$ws.Range("E13").Value = "  -4.11%  "
$ws.Range("D14").Value = "26.35"
$ws.Range("E14").Value = "  -3.57%  "
$ws.Range("D15").Value = "0.0000172"
$ws.Range("E15").Value = "  -5.69%  "
$ws.Range("D16").Value = "2.859.89"
$ws.Range("E16").Value = "  -2.35%  "
$ws.Range("D17").Value = "62.266.78"
$ws.Range("E17").Value = "  -2.31%  "
$ws.Range("D18").Value = "2.430.38"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").Value = "10.95"
$ws.Range("E19").Value = "  -5.05%  "
$ws.Range("D20").Value = "7.09"
$ws.Range("E20").Value = "  -4.25%  "
$ws.Range("D21").Value = "331.29"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "4.10"
$ws.Range("E22").Value = "  -3.09%  "
$ws.Range("D23").Value = "1.95"
$ws.Range("E23").Value = "  -8.82%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").Value = "65.70"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("D26").Value = "632.96"
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("D27").Value = "9.19"
$ws.Range("E27").Value = "  +1.13%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.554.56"
$ws.Range("E28").Value = "  -1.88%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0949"
$ws.Range("E29").Value = "  -10.37%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").Value = "1.42"
$ws.Range("E31").Value = "  -7.28%  "
$ws.Range("D32").Value = "8.02"
$ws.Range("E32").Value = "  -4.58%  "
$ws.Range("E33").Value = "  -1.87%  "
$ws.Range("D34").Value = "1.90"
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("D35").Value = "4.96"
$ws.Range("E35").Value = "  -6.24%  "
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("D37").Value = "1.44"
$ws.Range("E37").Value = "  -7.11%  "
$ws.Range("D38").Value = "0.375"
$ws.Range("E38").Value = "  -2.75%  "
$ws.Range("D39").Value = "149.17"
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("D40").Value = "18.34"
$ws.Range("E40").Value = "  -2.81%  "
$ws.Range("D41").Value = "5.23"
$ws.Range("E41").Value = "  -4.82%  "
$ws.Range("E42").Value = "  -4.68%  "
$ws.Range("E43").Value = "  +1.89%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "2.47"
$ws.Range("E45").Value = "  -10.53%  "
$ws.Range("D46").Value = "142.74"
$ws.Range("E46").Value = "  -5.34%  "
$ws.Range("D47").Value = "3.64"
$ws.Range("E47").Value = "  -3.82%  "
$ws.Range("D48").Value = "0.0521"
$ws.Range("E48").Value = "  -4.23%  "
$ws.Range("D49").Value = "0.596"
$ws.Range("E49").Value = "  -1.81%  "
$ws.Range("D50").Value = "19.62"
$ws.Range("E50").Value = "  -8.87%  "
$ws.Range("D51").Value = "0.0₆0233"
$ws.Range("E51").Value = "  +2.64%  "
